$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "auth-app" service was added to the Rolling release column; insert a
# row for it right after "audit" (before "auth-basic") and shift the rest of
# the table down by one row.
$ws.Rows("7:7").Insert()

# Helper to populate columns A:E of a given row in one go ($null leaves a
# cell untouched/blank, matching the sparse per-release columns in the sheet).
function Set-Row($r, [object[]]$vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($vals[$i] -ne $null) {
            $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
        }
    }
}

Set-Row 7 @($null,$null,$null,$null,"auth-app")
Set-Row 8 @("auth-basic","auth-basic","auth-basic","auth-basic","auth-basic")
Set-Row 9 @("auth-bearer","auth-bearer","auth-bearer","auth-bearer","auth-bearer")
Set-Row 10 @("auth-machine","auth-machine","auth-machine","auth-machine","auth-machine")
Set-Row 11 @("eventhistory","eventhistory","auth-service","auth-service","auth-service")
Set-Row 12 @($null,$null,"clientlog","clientlog","clientlog")
Set-Row 13 @($null,$null,$null,"collaboration","collaboration")
Set-Row 14 @($null,$null,"eventhistory","eventhistory","eventhistory")
Set-Row 15 @("frontend","frontend","frontend","frontend","frontend")
Set-Row 16 @("gateway","gateway","gateway","gateway","gateway")
Set-Row 17 @("graph","graph","graph","graph","graph")
Set-Row 18 @("groups","groups","groups","groups","groups")
Set-Row 19 @("idm","idm","idm","idm","idm")
Set-Row 20 @("idp","idp","idp","idp","idp")
Set-Row 21 @("invitations","invitations","invitations","invitations","invitations")
Set-Row 22 @("nats","nats","nats","nats","nats")
Set-Row 23 @("notifications","notifications","notifications","notifications","notifications")
Set-Row 24 @("ocdav","ocdav","ocdav","ocdav","ocdav")
Set-Row 25 @($null,$null,"ocm","ocm","ocm")
Set-Row 26 @("ocs","ocs","ocs","ocs","ocs")
Set-Row 27 @("policies","policies","policies","policies","policies")
Set-Row 28 @("postprocessing","postprocessing","postprocessing","postprocessing","postprocessing")
Set-Row 29 @("proxy","proxy","proxy","proxy","proxy")
Set-Row 30 @("search","search","search","search","search")
Set-Row 31 @("settings","settings","settings","settings","settings")
Set-Row 32 @("sharing","sharing","sharing","sharing","sharing")
Set-Row 33 @($null,$null,"sse","sse","sse")
Set-Row 34 @("store","store","store","store","store")
Set-Row 35 @("storage-publiclink","storage-publiclink","storage-publiclink","storage-publiclink","storage-publiclink")
Set-Row 36 @("storage-shares","storage-shares","storage-shares","storage-shares","storage-shares")
Set-Row 37 @("storage-system","storage-system","storage-system","storage-system","storage-system")
Set-Row 38 @("storage-users","storage-users","storage-users","storage-users","storage-users")
Set-Row 39 @("thumbnails","thumbnails","thumbnails","thumbnails","thumbnails")
Set-Row 40 @("userlog","userlog","userlog","userlog","userlog")
Set-Row 41 @("users","users","users","users","users")
Set-Row 42 @("web","web","web","web","web")
Set-Row 43 @("webdav","webdav","webdav","webdav","webdav")
Set-Row 44 @("webfinger","webfinger","webfinger","webfinger","webfinger")

# Restore the selection/view state recorded in the saved workbook.
$ws.Range("F34").Select()
